$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose target text is numeric-looking ("319.81", "13.00", ...) must be
# pre-formatted as Text so Excel keeps the exact string (incl. trailing zeros)
# instead of silently coercing the assignment to a floating-point number.

$ws.Range("D2").Value = "48.075.17"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "2.496.37"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.81"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("E6").Value = "  -3.29%  "
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("E9").Value = "  -4.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.77"
$ws.Range("E10").Value = "  -3.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.01"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0802"
$ws.Range("E12").Value = "  -2.07%  "
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("E14").Value = "  -2.21%  "
$ws.Range("D15").Value = "2.888.83"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "2.499.81"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.830"
$ws.Range("E17").Value = "  -3.49%  "
$ws.Range("D18").Value = "47.916.70"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("B19").Value = "ImmutableX"
$ws.Range("C19").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.98"
$ws.Range("E19").Value = "  +9.31%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.00"
$ws.Range("E20").Value = "  -2.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.64"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").Value = "0.0₃0931"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.06"
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "271.29"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.52"
$ws.Range("E25").Value = "  -2.34%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.76"
$ws.Range("E27").Value = "  -1.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.21"
$ws.Range("E28").Value = "  -4.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.72"
$ws.Range("E29").Value = "  -4.78%  "
$ws.Range("E30").Value = "  -3.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.74"
$ws.Range("E31").Value = "  -1.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.22"
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.07"
$ws.Range("E34").Value = "  -4.69%  "
$ws.Range("E35").Value = "  -2.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0772"
$ws.Range("E36").Value = "  -2.56%  "
$ws.Range("E37").Value = "  -2.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.58"
$ws.Range("E38").Value = "  -3.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.86"
$ws.Range("E39").Value = "  -4.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "121.98"
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.29"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.110"
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0303"
$ws.Range("E44").Value = "  +0.88%  "
$ws.Range("D45").Value = "1.998.50"
$ws.Range("E45").Value = "  -0.20%  "
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.88"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.91"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.16"
$ws.Range("E50").Value = "  -1.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.87"
$ws.Range("E51").Value = "  -2.20%  "
